$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing R column values (rows with plain value swaps) ---
$ws.Cells.Item(128, 18).Value = 0
$ws.Cells.Item(169, 18).Value = 0
$ws.Cells.Item(267, 18).Value = 0
$ws.Cells.Item(314, 18).Value = 0
$ws.Cells.Item(332, 18).Value = 0
$ws.Cells.Item(360, 18).Value = 1
$ws.Cells.Item(408, 18).Value = 2
$ws.Cells.Item(441, 18).Value = 1
$ws.Cells.Item(471, 18).Value = 1
$ws.Cells.Item(486, 18).Value = 0
$ws.Cells.Item(520, 18).Value = 0
$ws.Cells.Item(529, 18).Value = 2
$ws.Cells.Item(553, 18).Value = 2
$ws.Cells.Item(699, 18).Value = 0
$ws.Cells.Item(827, 18).Value = 1
$ws.Cells.Item(844, 18).Value = 2
$ws.Cells.Item(872, 18).Value = 2
$ws.Cells.Item(899, 18).Value = 0
$ws.Cells.Item(1049, 18).Value = 0
$ws.Cells.Item(1117, 18).Value = 0
$ws.Cells.Item(1183, 18).Value = 0
$ws.Cells.Item(1197, 18).Value = 1
$ws.Cells.Item(1205, 18).Value = 0
$ws.Cells.Item(1227, 18).Value = 1
$ws.Cells.Item(1265, 18).Value = 0

# --- Update R column for rows 1269-1317 (previously empty inlineStr, now numeric) ---
$ws.Cells.Item(1269, 18).Value = 1
$ws.Cells.Item(1270, 18).Value = 0
$ws.Cells.Item(1271, 18).Value = 0
$ws.Cells.Item(1272, 18).Value = 0
$ws.Cells.Item(1273, 18).Value = 0
$ws.Cells.Item(1274, 18).Value = 0
$ws.Cells.Item(1275, 18).Value = 0
$ws.Cells.Item(1276, 18).Value = 0
$ws.Cells.Item(1277, 18).Value = 0
$ws.Cells.Item(1278, 18).Value = 0
$ws.Cells.Item(1279, 18).Value = 0
$ws.Cells.Item(1280, 18).Value = 0
$ws.Cells.Item(1281, 18).Value = 0
$ws.Cells.Item(1282, 18).Value = 0
$ws.Cells.Item(1283, 18).Value = 0
$ws.Cells.Item(1284, 18).Value = 0
$ws.Cells.Item(1285, 18).Value = 0
$ws.Cells.Item(1286, 18).Value = 0
$ws.Cells.Item(1287, 18).Value = 0
$ws.Cells.Item(1288, 18).Value = 0
$ws.Cells.Item(1289, 18).Value = 0
$ws.Cells.Item(1290, 18).Value = 0
$ws.Cells.Item(1291, 18).Value = 0
$ws.Cells.Item(1292, 18).Value = 0
$ws.Cells.Item(1293, 18).Value = 0
$ws.Cells.Item(1294, 18).Value = 0
$ws.Cells.Item(1295, 18).Value = 0
$ws.Cells.Item(1296, 18).Value = 0
$ws.Cells.Item(1297, 18).Value = 2
$ws.Cells.Item(1298, 18).Value = 0
$ws.Cells.Item(1299, 18).Value = 0
$ws.Cells.Item(1300, 18).Value = 0
$ws.Cells.Item(1301, 18).Value = 0
$ws.Cells.Item(1302, 18).Value = 0
$ws.Cells.Item(1303, 18).Value = 0
$ws.Cells.Item(1304, 18).Value = 0
$ws.Cells.Item(1305, 18).Value = 0
$ws.Cells.Item(1306, 18).Value = 0
$ws.Cells.Item(1307, 18).Value = 0
$ws.Cells.Item(1308, 18).Value = 0
$ws.Cells.Item(1309, 18).Value = 0
$ws.Cells.Item(1310, 18).Value = 0
$ws.Cells.Item(1311, 18).Value = 0
$ws.Cells.Item(1312, 18).Value = 0
$ws.Cells.Item(1313, 18).Value = 0
$ws.Cells.Item(1314, 18).Value = 0
$ws.Cells.Item(1315, 18).Value = 0
$ws.Cells.Item(1316, 18).Value = 0
$ws.Cells.Item(1317, 18).Value = 0

# --- Append new rows 1318-1340 ---
# Row 1318
$ws.Cells.Item(1318, 1).Value = 45595
$ws.Cells.Item(1318, 2).Value = 2550
$ws.Cells.Item(1318, 3).Value = 2594.75
$ws.Cells.Item(1318, 4).Value = 2501
$ws.Cells.Item(1318, 5).Value = 2550
$ws.Cells.Item(1318, 6).Value = 2550
$ws.Cells.Item(1318, 7).Value = 562529
$ws.Cells.Item(1318, 8).Value = 2024
$ws.Cells.Item(1318, 9).Value = 10
$ws.Cells.Item(1318, 10).Value = 30
$ws.Cells.Item(1318, 11).Value = 0
$ws.Cells.Item(1318, 12).Value = 0
$ws.Cells.Item(1318, 13).Value = 0
$ws.Cells.Item(1318, 14).Value = 44
$ws.Cells.Item(1318, 15).Value = 0
$ws.Cells.Item(1318, 16).Value = 2
$ws.Cells.Item(1318, 17).Value = 0
$ws.Cells.Item(1318, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1319
$ws.Cells.Item(1319, 1).Value = 45596
$ws.Cells.Item(1319, 2).Value = 2550
$ws.Cells.Item(1319, 3).Value = 2555.75
$ws.Cells.Item(1319, 4).Value = 2488.5
$ws.Cells.Item(1319, 5).Value = 2500
$ws.Cells.Item(1319, 6).Value = 2500
$ws.Cells.Item(1319, 7).Value = 488935
$ws.Cells.Item(1319, 8).Value = 2024
$ws.Cells.Item(1319, 9).Value = 10
$ws.Cells.Item(1319, 10).Value = 31
$ws.Cells.Item(1319, 11).Value = 0
$ws.Cells.Item(1319, 12).Value = 0
$ws.Cells.Item(1319, 13).Value = 0
$ws.Cells.Item(1319, 14).Value = 44
$ws.Cells.Item(1319, 15).Value = 0
$ws.Cells.Item(1319, 16).Value = 0
$ws.Cells.Item(1319, 17).Value = 0
$ws.Cells.Item(1319, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1320
$ws.Cells.Item(1320, 1).Value = 45597
$ws.Cells.Item(1320, 2).Value = 2518
$ws.Cells.Item(1320, 3).Value = 2561.949951171875
$ws.Cells.Item(1320, 4).Value = 2511.75
$ws.Cells.Item(1320, 5).Value = 2551.14990234375
$ws.Cells.Item(1320, 6).Value = 2551.14990234375
$ws.Cells.Item(1320, 7).Value = 105769
$ws.Cells.Item(1320, 8).Value = 2024
$ws.Cells.Item(1320, 9).Value = 11
$ws.Cells.Item(1320, 10).Value = 1
$ws.Cells.Item(1320, 11).Value = 0
$ws.Cells.Item(1320, 12).Value = 0
$ws.Cells.Item(1320, 13).Value = 0
$ws.Cells.Item(1320, 14).Value = 44
$ws.Cells.Item(1320, 15).Value = 0
$ws.Cells.Item(1320, 16).Value = 0
$ws.Cells.Item(1320, 17).Value = 0
$ws.Cells.Item(1320, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1321
$ws.Cells.Item(1321, 1).Value = 45600
$ws.Cells.Item(1321, 2).Value = 2574.85009765625
$ws.Cells.Item(1321, 3).Value = 2578.85009765625
$ws.Cells.Item(1321, 4).Value = 2421
$ws.Cells.Item(1321, 5).Value = 2462.89990234375
$ws.Cells.Item(1321, 6).Value = 2462.89990234375
$ws.Cells.Item(1321, 7).Value = 535328
$ws.Cells.Item(1321, 8).Value = 2024
$ws.Cells.Item(1321, 9).Value = 11
$ws.Cells.Item(1321, 10).Value = 4
$ws.Cells.Item(1321, 11).Value = 0
$ws.Cells.Item(1321, 12).Value = 0
$ws.Cells.Item(1321, 13).Value = 0
$ws.Cells.Item(1321, 14).Value = 45
$ws.Cells.Item(1321, 15).Value = 0
$ws.Cells.Item(1321, 16).Value = 0
$ws.Cells.Item(1321, 17).Value = 0
$ws.Cells.Item(1321, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1322
$ws.Cells.Item(1322, 1).Value = 45601
$ws.Cells.Item(1322, 2).Value = 2460
$ws.Cells.Item(1322, 3).Value = 2472.14990234375
$ws.Cells.Item(1322, 4).Value = 2402
$ws.Cells.Item(1322, 5).Value = 2433.14990234375
$ws.Cells.Item(1322, 6).Value = 2433.14990234375
$ws.Cells.Item(1322, 7).Value = 279317
$ws.Cells.Item(1322, 8).Value = 2024
$ws.Cells.Item(1322, 9).Value = 11
$ws.Cells.Item(1322, 10).Value = 5
$ws.Cells.Item(1322, 11).Value = 0
$ws.Cells.Item(1322, 12).Value = 0
$ws.Cells.Item(1322, 13).Value = 0
$ws.Cells.Item(1322, 14).Value = 45
$ws.Cells.Item(1322, 15).Value = 0
$ws.Cells.Item(1322, 16).Value = 0
$ws.Cells.Item(1322, 17).Value = 0
$ws.Cells.Item(1322, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1323
$ws.Cells.Item(1323, 1).Value = 45602
$ws.Cells.Item(1323, 2).Value = 2465
$ws.Cells.Item(1323, 3).Value = 2473.550048828125
$ws.Cells.Item(1323, 4).Value = 2420.699951171875
$ws.Cells.Item(1323, 5).Value = 2462.800048828125
$ws.Cells.Item(1323, 6).Value = 2462.800048828125
$ws.Cells.Item(1323, 7).Value = 360009
$ws.Cells.Item(1323, 8).Value = 2024
$ws.Cells.Item(1323, 9).Value = 11
$ws.Cells.Item(1323, 10).Value = 6
$ws.Cells.Item(1323, 11).Value = 0
$ws.Cells.Item(1323, 12).Value = 0
$ws.Cells.Item(1323, 13).Value = 0
$ws.Cells.Item(1323, 14).Value = 45
$ws.Cells.Item(1323, 15).Value = 0
$ws.Cells.Item(1323, 16).Value = 0
$ws.Cells.Item(1323, 17).Value = 0
$ws.Cells.Item(1323, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1324
$ws.Cells.Item(1324, 1).Value = 45603
$ws.Cells.Item(1324, 2).Value = 2469.949951171875
$ws.Cells.Item(1324, 3).Value = 2473.35009765625
$ws.Cells.Item(1324, 4).Value = 2400
$ws.Cells.Item(1324, 5).Value = 2408.89990234375
$ws.Cells.Item(1324, 6).Value = 2408.89990234375
$ws.Cells.Item(1324, 7).Value = 236724
$ws.Cells.Item(1324, 8).Value = 2024
$ws.Cells.Item(1324, 9).Value = 11
$ws.Cells.Item(1324, 10).Value = 7
$ws.Cells.Item(1324, 11).Value = 0
$ws.Cells.Item(1324, 12).Value = 0
$ws.Cells.Item(1324, 13).Value = 0
$ws.Cells.Item(1324, 14).Value = 45
$ws.Cells.Item(1324, 15).Value = 0
$ws.Cells.Item(1324, 16).Value = 0
$ws.Cells.Item(1324, 17).Value = 0
$ws.Cells.Item(1324, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1325
$ws.Cells.Item(1325, 1).Value = 45604
$ws.Cells.Item(1325, 2).Value = 2413
$ws.Cells.Item(1325, 3).Value = 2425.449951171875
$ws.Cells.Item(1325, 4).Value = 2361.699951171875
$ws.Cells.Item(1325, 5).Value = 2373.449951171875
$ws.Cells.Item(1325, 6).Value = 2373.449951171875
$ws.Cells.Item(1325, 7).Value = 279067
$ws.Cells.Item(1325, 8).Value = 2024
$ws.Cells.Item(1325, 9).Value = 11
$ws.Cells.Item(1325, 10).Value = 8
$ws.Cells.Item(1325, 11).Value = 0
$ws.Cells.Item(1325, 12).Value = 0
$ws.Cells.Item(1325, 13).Value = 0
$ws.Cells.Item(1325, 14).Value = 45
$ws.Cells.Item(1325, 15).Value = 0
$ws.Cells.Item(1325, 16).Value = 0
$ws.Cells.Item(1325, 17).Value = 0
$ws.Cells.Item(1325, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1326
$ws.Cells.Item(1326, 1).Value = 45607
$ws.Cells.Item(1326, 2).Value = 2373.449951171875
$ws.Cells.Item(1326, 3).Value = 2394
$ws.Cells.Item(1326, 4).Value = 2351.800048828125
$ws.Cells.Item(1326, 5).Value = 2369.199951171875
$ws.Cells.Item(1326, 6).Value = 2369.199951171875
$ws.Cells.Item(1326, 7).Value = 100805
$ws.Cells.Item(1326, 8).Value = 2024
$ws.Cells.Item(1326, 9).Value = 11
$ws.Cells.Item(1326, 10).Value = 11
$ws.Cells.Item(1326, 11).Value = 0
$ws.Cells.Item(1326, 12).Value = 0
$ws.Cells.Item(1326, 13).Value = 0
$ws.Cells.Item(1326, 14).Value = 46
$ws.Cells.Item(1326, 15).Value = 0
$ws.Cells.Item(1326, 16).Value = 0
$ws.Cells.Item(1326, 17).Value = 0
$ws.Cells.Item(1326, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1327
$ws.Cells.Item(1327, 1).Value = 45608
$ws.Cells.Item(1327, 2).Value = 2379.949951171875
$ws.Cells.Item(1327, 3).Value = 2387.64990234375
$ws.Cells.Item(1327, 4).Value = 2300
$ws.Cells.Item(1327, 5).Value = 2308.14990234375
$ws.Cells.Item(1327, 6).Value = 2308.14990234375
$ws.Cells.Item(1327, 7).Value = 158104
$ws.Cells.Item(1327, 8).Value = 2024
$ws.Cells.Item(1327, 9).Value = 11
$ws.Cells.Item(1327, 10).Value = 12
$ws.Cells.Item(1327, 11).Value = 0
$ws.Cells.Item(1327, 12).Value = 0
$ws.Cells.Item(1327, 13).Value = 0
$ws.Cells.Item(1327, 14).Value = 46
$ws.Cells.Item(1327, 15).Value = 0
$ws.Cells.Item(1327, 16).Value = 0
$ws.Cells.Item(1327, 17).Value = 0
$ws.Cells.Item(1327, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1328
$ws.Cells.Item(1328, 1).Value = 45609
$ws.Cells.Item(1328, 2).Value = 2308.14990234375
$ws.Cells.Item(1328, 3).Value = 2332
$ws.Cells.Item(1328, 4).Value = 2283.199951171875
$ws.Cells.Item(1328, 5).Value = 2297.60009765625
$ws.Cells.Item(1328, 6).Value = 2297.60009765625
$ws.Cells.Item(1328, 7).Value = 286532
$ws.Cells.Item(1328, 8).Value = 2024
$ws.Cells.Item(1328, 9).Value = 11
$ws.Cells.Item(1328, 10).Value = 13
$ws.Cells.Item(1328, 11).Value = 0
$ws.Cells.Item(1328, 12).Value = 0
$ws.Cells.Item(1328, 13).Value = 0
$ws.Cells.Item(1328, 14).Value = 46
$ws.Cells.Item(1328, 15).Value = 0
$ws.Cells.Item(1328, 16).Value = 0
$ws.Cells.Item(1328, 17).Value = 0
$ws.Cells.Item(1328, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1329
$ws.Cells.Item(1329, 1).Value = 45610
$ws.Cells.Item(1329, 2).Value = 2307.550048828125
$ws.Cells.Item(1329, 3).Value = 2335.300048828125
$ws.Cells.Item(1329, 4).Value = 2276.449951171875
$ws.Cells.Item(1329, 5).Value = 2284.14990234375
$ws.Cells.Item(1329, 6).Value = 2284.14990234375
$ws.Cells.Item(1329, 7).Value = 250209
$ws.Cells.Item(1329, 8).Value = 2024
$ws.Cells.Item(1329, 9).Value = 11
$ws.Cells.Item(1329, 10).Value = 14
$ws.Cells.Item(1329, 11).Value = 0
$ws.Cells.Item(1329, 12).Value = 0
$ws.Cells.Item(1329, 13).Value = 0
$ws.Cells.Item(1329, 14).Value = 46
$ws.Cells.Item(1329, 15).Value = 0
$ws.Cells.Item(1329, 16).Value = 0
$ws.Cells.Item(1329, 17).Value = 0
$ws.Cells.Item(1329, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1330
$ws.Cells.Item(1330, 1).Value = 45614
$ws.Cells.Item(1330, 2).Value = 2280
$ws.Cells.Item(1330, 3).Value = 2288.800048828125
$ws.Cells.Item(1330, 4).Value = 2234.14990234375
$ws.Cells.Item(1330, 5).Value = 2242
$ws.Cells.Item(1330, 6).Value = 2242
$ws.Cells.Item(1330, 7).Value = 284801
$ws.Cells.Item(1330, 8).Value = 2024
$ws.Cells.Item(1330, 9).Value = 11
$ws.Cells.Item(1330, 10).Value = 18
$ws.Cells.Item(1330, 11).Value = 0
$ws.Cells.Item(1330, 12).Value = 0
$ws.Cells.Item(1330, 13).Value = 0
$ws.Cells.Item(1330, 14).Value = 47
$ws.Cells.Item(1330, 15).Value = 2
$ws.Cells.Item(1330, 16).Value = 0
$ws.Cells.Item(1330, 17).Value = 0
$ws.Cells.Item(1330, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1331
$ws.Cells.Item(1331, 1).Value = 45615
$ws.Cells.Item(1331, 2).Value = 2245.050048828125
$ws.Cells.Item(1331, 3).Value = 2312
$ws.Cells.Item(1331, 4).Value = 2245.050048828125
$ws.Cells.Item(1331, 5).Value = 2263.25
$ws.Cells.Item(1331, 6).Value = 2263.25
$ws.Cells.Item(1331, 7).Value = 227007
$ws.Cells.Item(1331, 8).Value = 2024
$ws.Cells.Item(1331, 9).Value = 11
$ws.Cells.Item(1331, 10).Value = 19
$ws.Cells.Item(1331, 11).Value = 0
$ws.Cells.Item(1331, 12).Value = 0
$ws.Cells.Item(1331, 13).Value = 0
$ws.Cells.Item(1331, 14).Value = 47
$ws.Cells.Item(1331, 15).Value = 0
$ws.Cells.Item(1331, 16).Value = 0
$ws.Cells.Item(1331, 17).Value = 0
$ws.Cells.Item(1331, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1332
$ws.Cells.Item(1332, 1).Value = 45617
$ws.Cells.Item(1332, 2).Value = 2251.050048828125
$ws.Cells.Item(1332, 3).Value = 2274.449951171875
$ws.Cells.Item(1332, 4).Value = 2244
$ws.Cells.Item(1332, 5).Value = 2253.949951171875
$ws.Cells.Item(1332, 6).Value = 2253.949951171875
$ws.Cells.Item(1332, 7).Value = 202738
$ws.Cells.Item(1332, 8).Value = 2024
$ws.Cells.Item(1332, 9).Value = 11
$ws.Cells.Item(1332, 10).Value = 21
$ws.Cells.Item(1332, 11).Value = 0
$ws.Cells.Item(1332, 12).Value = 0
$ws.Cells.Item(1332, 13).Value = 0
$ws.Cells.Item(1332, 14).Value = 47
$ws.Cells.Item(1332, 15).Value = 0
$ws.Cells.Item(1332, 16).Value = 0
$ws.Cells.Item(1332, 17).Value = 0
$ws.Cells.Item(1332, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1333
$ws.Cells.Item(1333, 1).Value = 45618
$ws.Cells.Item(1333, 2).Value = 2253.949951171875
$ws.Cells.Item(1333, 3).Value = 2287.60009765625
$ws.Cells.Item(1333, 4).Value = 2244.10009765625
$ws.Cells.Item(1333, 5).Value = 2278.85009765625
$ws.Cells.Item(1333, 6).Value = 2278.85009765625
$ws.Cells.Item(1333, 7).Value = 180507
$ws.Cells.Item(1333, 8).Value = 2024
$ws.Cells.Item(1333, 9).Value = 11
$ws.Cells.Item(1333, 10).Value = 22
$ws.Cells.Item(1333, 11).Value = 0
$ws.Cells.Item(1333, 12).Value = 0
$ws.Cells.Item(1333, 13).Value = 0
$ws.Cells.Item(1333, 14).Value = 47
$ws.Cells.Item(1333, 15).Value = 0
$ws.Cells.Item(1333, 16).Value = 0
$ws.Cells.Item(1333, 17).Value = 0
$ws.Cells.Item(1333, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1334
$ws.Cells.Item(1334, 1).Value = 45621
$ws.Cells.Item(1334, 2).Value = 2300.949951171875
$ws.Cells.Item(1334, 3).Value = 2320.949951171875
$ws.Cells.Item(1334, 4).Value = 2255.75
$ws.Cells.Item(1334, 5).Value = 2309.39990234375
$ws.Cells.Item(1334, 6).Value = 2309.39990234375
$ws.Cells.Item(1334, 7).Value = 691659
$ws.Cells.Item(1334, 8).Value = 2024
$ws.Cells.Item(1334, 9).Value = 11
$ws.Cells.Item(1334, 10).Value = 25
$ws.Cells.Item(1334, 11).Value = 0
$ws.Cells.Item(1334, 12).Value = 0
$ws.Cells.Item(1334, 13).Value = 0
$ws.Cells.Item(1334, 14).Value = 48
$ws.Cells.Item(1334, 15).Value = 0
$ws.Cells.Item(1334, 16).Value = 0
$ws.Cells.Item(1334, 17).Value = 0
$ws.Cells.Item(1334, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1335
$ws.Cells.Item(1335, 1).Value = 45622
$ws.Cells.Item(1335, 2).Value = 2309.699951171875
$ws.Cells.Item(1335, 3).Value = 2355.64990234375
$ws.Cells.Item(1335, 4).Value = 2309.699951171875
$ws.Cells.Item(1335, 5).Value = 2334.14990234375
$ws.Cells.Item(1335, 6).Value = 2334.14990234375
$ws.Cells.Item(1335, 7).Value = 380827
$ws.Cells.Item(1335, 8).Value = 2024
$ws.Cells.Item(1335, 9).Value = 11
$ws.Cells.Item(1335, 10).Value = 26
$ws.Cells.Item(1335, 11).Value = 0
$ws.Cells.Item(1335, 12).Value = 0
$ws.Cells.Item(1335, 13).Value = 0
$ws.Cells.Item(1335, 14).Value = 48
$ws.Cells.Item(1335, 15).Value = 0
$ws.Cells.Item(1335, 16).Value = 0
$ws.Cells.Item(1335, 17).Value = 0
$ws.Cells.Item(1335, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1336
$ws.Cells.Item(1336, 1).Value = 45623
$ws.Cells.Item(1336, 2).Value = 2340.39990234375
$ws.Cells.Item(1336, 3).Value = 2380.949951171875
$ws.Cells.Item(1336, 4).Value = 2324.800048828125
$ws.Cells.Item(1336, 5).Value = 2367.14990234375
$ws.Cells.Item(1336, 6).Value = 2367.14990234375
$ws.Cells.Item(1336, 7).Value = 190150
$ws.Cells.Item(1336, 8).Value = 2024
$ws.Cells.Item(1336, 9).Value = 11
$ws.Cells.Item(1336, 10).Value = 27
$ws.Cells.Item(1336, 11).Value = 0
$ws.Cells.Item(1336, 12).Value = 0
$ws.Cells.Item(1336, 13).Value = 0
$ws.Cells.Item(1336, 14).Value = 48
$ws.Cells.Item(1336, 15).Value = 0
$ws.Cells.Item(1336, 16).Value = 0
$ws.Cells.Item(1336, 17).Value = 0
$ws.Cells.Item(1336, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1337
$ws.Cells.Item(1337, 1).Value = 45624
$ws.Cells.Item(1337, 2).Value = 2389
$ws.Cells.Item(1337, 3).Value = 2394.60009765625
$ws.Cells.Item(1337, 4).Value = 2335
$ws.Cells.Item(1337, 5).Value = 2341.39990234375
$ws.Cells.Item(1337, 6).Value = 2341.39990234375
$ws.Cells.Item(1337, 7).Value = 268139
$ws.Cells.Item(1337, 8).Value = 2024
$ws.Cells.Item(1337, 9).Value = 11
$ws.Cells.Item(1337, 10).Value = 28
$ws.Cells.Item(1337, 11).Value = 0
$ws.Cells.Item(1337, 12).Value = 0
$ws.Cells.Item(1337, 13).Value = 0
$ws.Cells.Item(1337, 14).Value = 48
$ws.Cells.Item(1337, 15).Value = 0
$ws.Cells.Item(1337, 16).Value = 0
$ws.Cells.Item(1337, 17).Value = 0
$ws.Cells.Item(1337, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1338
$ws.Cells.Item(1338, 1).Value = 45625
$ws.Cells.Item(1338, 2).Value = 2344.699951171875
$ws.Cells.Item(1338, 3).Value = 2355.89990234375
$ws.Cells.Item(1338, 4).Value = 2323.050048828125
$ws.Cells.Item(1338, 5).Value = 2344.89990234375
$ws.Cells.Item(1338, 6).Value = 2344.89990234375
$ws.Cells.Item(1338, 7).Value = 98585
$ws.Cells.Item(1338, 8).Value = 2024
$ws.Cells.Item(1338, 9).Value = 11
$ws.Cells.Item(1338, 10).Value = 29
$ws.Cells.Item(1338, 11).Value = 0
$ws.Cells.Item(1338, 12).Value = 0
$ws.Cells.Item(1338, 13).Value = 0
$ws.Cells.Item(1338, 14).Value = 48
$ws.Cells.Item(1338, 15).Value = 0
$ws.Cells.Item(1338, 16).Value = 0
$ws.Cells.Item(1338, 17).Value = 0
$ws.Cells.Item(1338, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1339
$ws.Cells.Item(1339, 1).Value = 45628
$ws.Cells.Item(1339, 2).Value = 2344.89990234375
$ws.Cells.Item(1339, 3).Value = 2449
$ws.Cells.Item(1339, 4).Value = 2251
$ws.Cells.Item(1339, 5).Value = 2362.800048828125
$ws.Cells.Item(1339, 6).Value = 2362.800048828125
$ws.Cells.Item(1339, 7).Value = 174198
$ws.Cells.Item(1339, 8).Value = 2024
$ws.Cells.Item(1339, 9).Value = 12
$ws.Cells.Item(1339, 10).Value = 2
$ws.Cells.Item(1339, 11).Value = 0
$ws.Cells.Item(1339, 12).Value = 0
$ws.Cells.Item(1339, 13).Value = 0
$ws.Cells.Item(1339, 14).Value = 49
$ws.Cells.Item(1339, 15).Value = 0
$ws.Cells.Item(1339, 16).Value = 0
$ws.Cells.Item(1339, 17).Value = 0
$ws.Cells.Item(1339, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1340
$ws.Cells.Item(1340, 1).Value = 45629
$ws.Cells.Item(1340, 2).Value = 2356
$ws.Cells.Item(1340, 3).Value = 2387.699951171875
$ws.Cells.Item(1340, 4).Value = 2339.10009765625
$ws.Cells.Item(1340, 5).Value = 2357.550048828125
$ws.Cells.Item(1340, 6).Value = 2357.550048828125
$ws.Cells.Item(1340, 7).Value = 143183
$ws.Cells.Item(1340, 8).Value = 2024
$ws.Cells.Item(1340, 9).Value = 12
$ws.Cells.Item(1340, 10).Value = 3
$ws.Cells.Item(1340, 11).Value = 0
$ws.Cells.Item(1340, 12).Value = 0
$ws.Cells.Item(1340, 13).Value = 0
$ws.Cells.Item(1340, 14).Value = 49
$ws.Cells.Item(1340, 15).Value = 0
$ws.Cells.Item(1340, 16).Value = 0
$ws.Cells.Item(1340, 17).Value = 0
$ws.Cells.Item(1340, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

